$d = $word.ActiveDocument
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

# ------------------------------------------------------------------
# Change 1: expand the placeholder "Текст завдання" paragraph into the
# full assignment-description paragraph (with No/Exception spell-check
# markers).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Текст завдання`r") {
        $xml1 = "<w:p $ns w14:paraId='1DE3BF90' w14:textId='3C4D65ED' w:rsidR='008E3F95' w:rsidRPr='0011071F' w:rsidRDefault='008E3F95' w:rsidP='0011071F'>" +
                "<w:pPr><w:spacing w:after='0' w:line='360' w:lineRule='auto'/><w:ind w:firstLine='851'/></w:pPr>" +
                "<w:r w:rsidRPr='0011071F'><w:t xml:space='preserve'>Для  завдання  з  лабораторної  роботи  </w:t></w:r>" +
                "<w:proofErr w:type='spellStart'/>" +
                "<w:r w:rsidRPr='0011071F'><w:t>No</w:t></w:r>" +
                "<w:proofErr w:type='spellEnd'/>" +
                "<w:r w:rsidRPr='0011071F'><w:t xml:space='preserve'>  5  виконати  обробку виняткових ситуацій з використанням класу </w:t></w:r>" +
                "<w:proofErr w:type='spellStart'/>" +
                "<w:r w:rsidRPr='0011071F'><w:t>Exception</w:t></w:r>" +
                "<w:proofErr w:type='spellEnd'/>" +
                "<w:r w:rsidRPr='0011071F'><w:t>.</w:t></w:r>" +
                "</w:p>"
        $para.Range.InsertXML($xml1)
        break
    }
}

# ------------------------------------------------------------------
# Change 2: mark the "Висновки" paragraph's first run with a
# lastRenderedPageBreak (page break landed there during pagination).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Таким чином, ми навчилися")) {
        $xml2 = "<w:p $ns w14:paraId='2CFD46BA' w14:textId='042F0019' w:rsidR='008E3F95' w:rsidRPr='0011071F' w:rsidRDefault='008E3F95' w:rsidP='0011071F'>" +
                "<w:pPr><w:spacing w:after='0' w:line='360' w:lineRule='auto'/><w:ind w:firstLine='851'/></w:pPr>" +
                "<w:r w:rsidRPr='0011071F'><w:lastRenderedPageBreak/><w:t xml:space='preserve'>Таким чином, ми навчилися </w:t></w:r>" +
                "<w:proofErr w:type='spellStart'/>" +
                "<w:r w:rsidR='005F649F' w:rsidRPr='00B13FB2'><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>обробляти</w:t></w:r>" +
                "<w:proofErr w:type='spellEnd'/>" +
                "<w:r w:rsidR='005F649F' w:rsidRPr='00B13FB2'><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r>" +
                "<w:proofErr w:type='spellStart'/>" +
                "<w:r w:rsidR='005F649F' w:rsidRPr='00B13FB2'><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>виняткові</w:t></w:r>" +
                "<w:proofErr w:type='spellEnd'/>" +
                "<w:r w:rsidR='005F649F' w:rsidRPr='00B13FB2'><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r>" +
                "<w:proofErr w:type='spellStart'/>" +
                "<w:r w:rsidR='005F649F' w:rsidRPr='00B13FB2'><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>ситуації</w:t></w:r>" +
                "<w:proofErr w:type='spellEnd'/>" +
                "<w:r w:rsidR='005F649F' w:rsidRPr='0011071F'><w:t xml:space='preserve'> </w:t></w:r>" +
                "<w:r w:rsidRPr='0011071F'><w:t>при створенні програм мовою програмування С++.</w:t></w:r>" +
                "</w:p>"
        $para.Range.InsertXML($xml2)
        break
    }
}
